# "alterado novamente slide 2"
#
# 1) Slide 2 title/body text is updated.
# 2) The cached caption of the slide-number field placeholders (shown as
#    <nbr> on the slide master and every slide layout) is refreshed from
#    the old glyph to the new one.

$p = $ppt.ActivePresentation

$oldCaption = [string]([char]0x2039) + "n" + [char]0xba + [char]0x203a
$newCaption = [string]([char]0x2039) + "#" + [char]0x203a

# --- helper: update the slide-number placeholder caption on a Shapes
#     collection (slide master or a custom layout) -------------------
function Update-SlideNumberCaption {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($shp.PlaceholderFormat.Type -eq 13 -and $tr.Text -eq $oldCaption) {
                $tr.Text = $newCaption
            }
        }
    }
}

# Slide master
Update-SlideNumberCaption $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-SlideNumberCaption $layouts.Item($li).Shapes
}

# --- Slide 2 text edits -------------------------------------------------
$s2 = $p.Slides.Item(2)

$title = $s2.Shapes.Item(1).TextFrame.TextRange
$title.Runs(1).Text = "Alterado novamente slide "
$title.Runs(2).Text = "2"

$body = $s2.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Alterado novamente"
